$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) holds a date serial value that needs to change
# from 45186 (2023-09-17) to 45188 (2023-09-19) for every data row,
# from row 2 through row 385.
$ws.Range("C2:C385").Value = 45188
